# Cotações atualizadas - 2025-10-19
# Appends a new row (45) to the quotes table with the values for 2025-10-19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 45

# Column A: date serial 45949 (2025-10-19), same number format/style as the
# other date cells in column A (style index "2" -> YYYY-MM-DD HH:MM:SS).
$ws.Cells.Item($newRow, 1).Value = 45949
$ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($newRow - 1, 1).NumberFormat

# Columns B-E: quote values stored as plain text (Portuguese decimal comma),
# matching the previous row's values.
$ws.Cells.Item($newRow, 2).Value = "21,7414"
$ws.Cells.Item($newRow, 3).Value = "15,4996"
$ws.Cells.Item($newRow, 4).Value = "15,5859"
$ws.Cells.Item($newRow, 5).Value = "15,5859"
